$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.013.65"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.742.04"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'246.73"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "'0.5049"
$ws.Range("E7").Value = "  -4.26%  "

$ws.Range("D8").Value = "'0.2742"
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").Value = "1.749.55"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").Value = "'0.07245"
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").Value = "'0.6521"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "'4.680"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").Value = "'77.56"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "26.023.31"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "'11.90"
$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").Value = "'0.000006867"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").Value = "1.972.25"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'4.485"
$ws.Range("E22").Value = "  +2.89%  "

$ws.Range("D23").Value = "'8.712"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("D24").Value = "'5.372"
$ws.Range("E24").Value = "  +2.28%  "

$ws.Range("D25").Value = "'135.89"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "'1.500"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").Value = "'15.27"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "'1.783"
$ws.Range("E28").Value = "  +1.23%  "

$ws.Range("D29").Value = "'105.33"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "'3.894"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").Value = "'0.08180"
$ws.Range("E31").Value = "  -3.03%  "

$ws.Range("D32").Value = "'3.655"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("D33").Value = "'0.04670"
$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("D34").Value = "'2.657"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'0.9981"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6103"
$ws.Range("E36").Value = "  -2.45%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.772"
$ws.Range("E37").Value = "  +2.74%  "

$ws.Range("D38").Value = "'0.01623"
$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("D39").Value = "'1.930"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "'100.72"
$ws.Range("E41").Value = "  +1.69%  "

$ws.Range("D42").Value = "'0.3926"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("D43").Value = "'0.7638"
$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("D44").Value = "'5.004"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").Value = "'0.1165"
$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("D46").Value = "'6.331"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").Value = "'55.56"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("D48").Value = "'0.05301"
$ws.Range("E48").Value = "  -0.41%  "

$ws.Range("D49").Value = "'30.65"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").Value = "'0.3466"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "'7.575"
$ws.Range("E51").Value = "  +0.40%  "
